# Update column F (dSF) values on the active worksheet to reflect the
# repulled/recomputed data, per the commit "repull data, push all data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 2
    4  = 4
    5  = 1
    6  = 8
    7  = 6
    9  = -2
    10 = -1
    12 = 1
    14 = -2
    15 = -1
    16 = 6
    17 = 1
    18 = 4
    19 = -2
    20 = -5
    21 = -5
    23 = -3
    25 = -5
    26 = 2
    27 = -4
    28 = -3
    29 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
